# ---------------------------------------------------------------------------
# Applies the commit "Add files via upload / пофиксил описание":
#
#   1. Reword the sentence describing how the property-based solution was
#      chosen ("решить это ... тут решено" -> "сделать реализовать логику
#      ... решено").
#   2. Relocate the auto-maintained hidden "_GoBack" bookmark from its old
#      spot (mid-word, splitting "сво|йства") to its new spot (mid-word,
#      splitting "касто|много" inside "кастомного" further down). Moving
#      this bookmark is what produced the large, otherwise textually-inert,
#      run churn in the raw XML diff: wherever "_GoBack" sits, Word has to
#      break the enclosing run in two around it, and merges the runs back
#      together once the bookmark leaves.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Reword the solution-choice sentence --------------------------------
$d.Content.Find.Execute(
    "решить это через свойство продукта в корзине, тут",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "сделать реализовать логику, через свойство продукта в корзине,",
    2
) | Out-Null

# --- 2a. Drop the old hidden bookmark, wherever it currently sits ----------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2b. Re-join the run the old bookmark used to split ("сво" | "йства") --
# Re-typing the whole phrase in one Find/Replace pass makes Word emit it as
# a single run again (the "EVERY" that follows keeps its own run/formatting
# because the replacement text stops right before it).
$d.Content.Find.Execute(
    "которые будут подходить под наши условия, устанавливать свойства ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "которые будут подходить под наши условия, устанавливать свойства ",
    2
) | Out-Null

# --- 2c. Find the new split point, inside the second "кастомного" ----------
$rng = $d.Content
$rng.Find.Execute("кастомного", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Find.Execute("кастомного", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $rng.Start + 5   # right after "касто", right before "много"

# --- 2d. Re-create "_GoBack", collapsed, at the new location ---------------
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
